$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI statistics (E/K expressing-cell counts changed 1 -> 3,
# with derived average/total/specificity/edge-weight columns recomputed).
$data = @(
    ,@(2,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,6.278248333333334,18.834745,0.02606887289604514,0.02606887289604514,76.62885867658001,689.65972808922,0.005261475361484037,0.005261475361484037)
    ,@(3,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,8.908825666666667,26.726477,0.0369916944387659,0.03699169443876589,108.736244050868,978.6261964578118,0.007466026231561395,0.007466026231561392)
    ,@(4,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,64.156346,192.469038,0.2663933537749559,0.2663933537749558,783.057201598392,7047.514814385527,0.05376611689117824,0.05376611689117822)
    ,@(5,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,60.88514866666666,182.655446,0.2528105161787131,0.2528105161787131,743.1307595638639,6688.176836074775,0.05102469551724103,0.05102469551724102)
    ,@(6,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,93.12082700000001,279.362481,0.3866611949942951,0.3866611949942951,1136.581784148804,10229.23605733923,0.07803975104014163,0.07803975104014163)
    ,@(7,3,12.205452,36.616356,0.2018297984138105,0.2018297984138104,3,7.483737333333333,22.451212,0.03107436771722491,0.03107436771722491,91.34239680260799,822.0815712234719,0.006271733372204123,0.006271733372204122)
    ,@(8,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,6.278248333333334,18.834745,0.02606887289604514,0.02606887289604514,273.5558154525056,2462.002339072551,0.01878283466374955,0.01878283466374955)
    ,@(9,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,8.908825666666667,26.726477,0.0369916944387659,0.03699169443876589,388.1753222519145,3493.57790026723,0.02665281630494627,0.02665281630494626)
    ,@(10,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,64.156346,192.469038,0.2663933537749559,0.2663933537749558,2795.420094057514,25158.78084651762,0.1919385751516641,0.191938575151664)
    ,@(11,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,60.88514866666666,182.655446,0.2528105161787131,0.2528105161787131,2652.887494753504,23875.98745278154,0.1821520303381561,0.1821520303381561)
    ,@(12,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,93.12082700000001,279.362481,0.3866611949942951,0.3866611949942951,4057.460363641244,36517.14327277119,0.2785925316152608,0.2785925316152608)
    ,@(13,3,43.57199666666667,130.71599,0.7205081224128821,0.720508122412882,3,7.483737333333333,22.451212,0.03107436771722491,0.03107436771722491,326.0813781422089,2934.73240327988,0.0223893343391052,0.02238933433910519)
    ,@(14,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,6.278248333333334,18.834745,0.02606887289604514,0.02606887289604514,29.48601512894111,265.37413616047,0.002024562870811548,0.002024562870811548)
    ,@(15,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,8.908825666666667,26.726477,0.0369916944387659,0.03699169443876589,41.84061452200689,376.565530698062,0.002872851902258237,0.002872851902258236)
    ,@(16,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,64.156346,192.469038,0.2663933537749559,0.2663933537749558,301.3125458465586,2711.812912619028,0.02068866173211354,0.02068866173211354)
    ,@(17,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,60.88514866666666,182.655446,0.2528105161787131,0.2528105161787131,285.9492519882529,2573.543267894276,0.01963379032331596,0.01963379032331596)
    ,@(18,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,93.12082700000001,279.362481,0.3866611949942951,0.3866611949942951,437.3452542747207,3936.107288472486,0.03002891233889264,0.03002891233889264)
    ,@(19,3,4.696535333333333,14.089606,0.07766207917330753,0.07766207917330753,3,7.483737333333333,22.451212,0.03107436771722491,0.03107436771722491,35.14763681138577,316.328731302472,0.002413300005915592,0.002413300005915592)
)

$cols = @('E','G','H','I','J','K','M','N','O','P','Q','R','S','T')

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).Value2 = $row[$i + 1]
    }
}

